$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The naive forecaster rolling window advanced by one period: the old
# row 2 (earliest forecast date) is dropped and every subsequent row
# shifts up by one, with C/E recomputed by the (bugfixed) forecaster
# module.
$ws.Rows.Item(2).Delete()

$data = @"
2,39583,2008,,2009,1.516248937663556
3,39765,2008,,2009,1.560682679516057
4,39948,2009,1.514319819128396,2010,1.602279001294704
5,40130,2009,1.834695583582491,2010,2.1453644888767
6,40310,2010,2.007652128026982,2011,1.768040115052738
7,40494,2010,1.767835936772166,2011,1.317672174811868
8,40676,2011,1.022680528298392,2012,1.510468690286459
9,40862,2011,1.074400434091016,2012,1.501816644427989
10,41044,2012,0.9070039918702477,2013,1.042579621507111
11,41228,2012,0.9212998022035679,2013,1.028888107831327
12,41409,2013,1.022042907336096,2014,1.083482333436536
13,41592,2013,1.141837882844188,2014,1.303605130836716
14,41774,2014,1.34489417553354,2015,1.095916825800991
15,41957,2014,1.335361538769475,2015,1.192378712846454
16,42137,2015,1.277042522796856,2016,1.293136192195643
17,42321,2015,1.202048372526998,2016,1.210961441871872
18,42503,2016,2.615369162917314,2017,1.825134644920934
19,42689,2016,2.677488680362305,2017,2.033218171624651
20,42867,2017,2.618053282882693,2018,2.406099663413808
21,43053,2017,2.466954516646402,2018,2.152537330144288
22,43145,2018,1.731723847815725,2019,2.356276715023498
23,43235,2018,1.431088640641853,2019,2.21629047761287
24,43326,2018,1.372961566907027,2019,1.902399534782662
25,43418,2018,1.401189216021326,2019,1.966855307908655
26,43510,2019,1.983559881711905,2020,1.93103453922987
27,43600,2019,2.136062314641141,2020,2.031764787322499
28,43691,2019,2.241561867365394,2020,2.135927826705641
29,43783,2019,2.217567799050979,2020,1.950353221540246
30,43875,2020,2.041276490941102,2021,1.854752869950294
31,43966,2020,2.166968775134936,2021,1.984987808509886
32,44068,2020,2.139672475020404,2021,2.011395609719546
33,44159,2020,2.139672475020404,2021,2.210985773414453
34,44251,2021,2.330672672271739,2022,2.307457288603798
35,44341,2021,2.459440348120401,2022,2.526389380645511
36,44432,2021,2.100991693542231,2022,1.091147151778871
37,44525,2021,2.100991693542231,2022,1.114171399050901
38,44617,2022,0.7608230790701942,2023,1.459415358104388
39,44706,2022,0.8171929556848756,2023,1.509741350988136
40,44798,2022,0.8967077601845341,2023,1.773412413757813
41,44890,2022,0.8967077601845341,2023,0.1338254721205745
42,44981,2023,-0.4883557973630492,2024,0.4753196237801127
43,45071,2023,0.8766015904249524,2024,2.477445663648559
44,45163,2023,0.782207885866093,2024,2.01691766737
45,45254,2023,0.782207885866093,2024,1.823564868738359
46,45345,2024,1.617535832906758,2025,1.142484412546874
47,45436,2024,1.554086551645839,2025,0.9888012784191602
48,45534,2024,1.508385007449875,2025,0.8780954941978392
49,45618,2024,1.508385007449875,2025,0.6266145540918089
50,45713,2025,0.3854686824285025,2026,0.6767639290315763
51,45800,2025,0.5837948599211717,2026,1.328924132093245
52,45891,2025,0.6014263374495288,2026,1.609787824259601
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
  $line = $line.Trim()
  if ($line -eq "") { continue }
  $parts = $line -split ","
  $r = [int]$parts[0]

  $ws.Cells.Item($r, 1).Value2 = [double]$parts[1]
  $ws.Cells.Item($r, 2).Value2 = [double]$parts[2]

  if ($parts[3] -ne "") {
    $ws.Cells.Item($r, 3).Value2 = [double]$parts[3]
  } else {
    $ws.Cells.Item($r, 3).ClearContents()
  }

  $ws.Cells.Item($r, 4).Value2 = [double]$parts[4]

  if ($parts[5] -ne "") {
    $ws.Cells.Item($r, 5).Value2 = [double]$parts[5]
  } else {
    $ws.Cells.Item($r, 5).ClearContents()
  }
}
